$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.890.90'
$ws.Range("E2").Value = '  -0.18%  '

$ws.Range("D3").Value = '2.975.07'
$ws.Range("E3").Value = '  +1.43%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '353.85'
$ws.Range("E5").Value = '  -1.04%  '

$ws.Range("D6").Value = '106.72'
$ws.Range("E6").Value = '  -3.68%  '

$ws.Range("E7").Value = '  -2.88%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").Value = '0.602'
$ws.Range("E9").Value = '  -5.37%  '

$ws.Range("D10").Value = '37.68'
$ws.Range("E10").Value = '  -4.66%  '

$ws.Range("E11").Value = '  +2.52%  '

$ws.Range("D12").Value = '0.0850'
$ws.Range("E12").Value = '  -3.66%  '

$ws.Range("D13").Value = '18.92'
$ws.Range("E13").Value = '  -4.32%  '

$ws.Range("D14").Value = '3.447.24'
$ws.Range("E14").Value = '  +1.54%  '

$ws.Range("E15").Value = '  -5.35%  '

$ws.Range("D16").Value = '2.962.82'
$ws.Range("E16").Value = '  +1.56%  '

$ws.Range("D17").Value = '0.987'
$ws.Range("E17").Value = '  -0.25%  '

$ws.Range("D18").Value = '51.827.88'
$ws.Range("E18").Value = '  -0.42%  '

$ws.Range("E19").Value = '  +0.66%  '

$ws.Range("E20").Value = '  -3.24%  '

$ws.Range("D21").Value = '13.38'
$ws.Range("E21").Value = '  -5.20%  '

$ws.Range("D22").Value = '0.0₃0963'
$ws.Range("E22").Value = '  -2.33%  '

$ws.Range("D23").Value = '68.93'
$ws.Range("E23").Value = '  -3.02%  '

$ws.Range("D24").Value = '262.83'
$ws.Range("E24").Value = '  -3.15%  '

$ws.Range("D25").Value = '2.70'
$ws.Range("E25").Value = '  -4.81%  '

$ws.Range("D26").Value = '0.177'
$ws.Range("E26").Value = '  -3.92%  '

$ws.Range("D27").Value = '26.71'
$ws.Range("E27").Value = '  -2.02%  '

$ws.Range("E28").Value = '  +0.01%  '

$ws.Range("B29").Value = 'Hedera'
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D29").Value = '0.111'
$ws.Range("E29").Value = '  +3.08%  '

$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").Value = '7.29'
$ws.Range("E30").Value = '  -3.73%  '

$ws.Range("D31").Value = '6.24'
$ws.Range("E31").Value = '  +2.51%  '

$ws.Range("D32").Value = '10.07'
$ws.Range("E32").Value = '  -5.54%  '

$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D33").Value = '35.68'
$ws.Range("E33").Value = '  -7.04%  '

$ws.Range("B34").Value = 'Toncoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D34").Value = '2.16'
$ws.Range("E34").Value = '  +12.45%  '

$ws.Range("D35").Value = '50.90'
$ws.Range("E35").Value = '  -2.83%  '

$ws.Range("E36").Value = '  -4.40%  '

$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  -0.12%  '

$ws.Range("E38").Value = '  -1.08%  '

$ws.Range("D39").Value = '2.83'
$ws.Range("E39").Value = '  +1.74%  '

$ws.Range("D40").Value = '1.93'
$ws.Range("E40").Value = '  -4.68%  '

$ws.Range("D41").Value = '17.35'
$ws.Range("E41").Value = '  -6.42%  '

$ws.Range("E42").Value = '  -3.69%  '

$ws.Range("D43").Value = '22.90'
$ws.Range("E43").Value = '  -0.98%  '

$ws.Range("D44").Value = '123.64'
$ws.Range("E44").Value = '  +3.70%  '

$ws.Range("E45").Value = '  -0.41%  '

$ws.Range("D46").Value = '2.102.86'
$ws.Range("E46").Value = '  -1.93%  '

$ws.Range("E47").Value = '  -5.96%  '

$ws.Range("D48").Value = '2.31'
$ws.Range("E48").Value = '  -7.89%  '

$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '3.272.36'
$ws.Range("E49").Value = '  +1.59%  '

$ws.Range("B50").Value = 'TheGraph'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D50").Value = '0.238'
$ws.Range("E50").Value = '  -4.01%  '

$ws.Range("B51").Value = 'BEAM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range("D51").Value = '0.0327'
$ws.Range("E51").Value = '  -2.25%  '
